$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Metadata sheet ("Metadata") ---
# Version bump: 5.0.0 -> 6.0.0
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Date refresh
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# Remove the duplicate "Contact / No display for ContactDetail" row (row 11),
# leaving row 10 which we then turn into "Jurisdiction / United States of America".
$ws1.Rows.Item(11).Delete()

$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"
